$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version: 5.0.0 -> 6.0.0
$ws.Range("B3").Value = "6.0.0"

# Date: 2021-12-16T17:36:56+00:00 -> 2022-01-21T20:46:54+00:00
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was blank -> "Alvearie Team"
$ws.Range("B9").Value = "Alvearie Team"

# Row 10 was a duplicated "Contact" / "No display for ContactDetail" row;
# turn it into "Jurisdiction" / "United States of America"
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# Row 11 was the second duplicated "Contact" row - remove it entirely,
# shifting the remaining rows up by one
$ws.Rows.Item(11).Delete()

# After the shift, the "Case Sensitive" row (now row 14) gets a value of "true".
# A plain Value assignment of the literal "true" gets auto-coerced to a Boolean
# by the engine, so instead compute it as a text formula in a scratch cell and
# paste-special just the value back in, which keeps it a real text string.
$helper = $ws.Range("ZZ1")
$helper.Formula = "=""true"""
$helper.Copy()
$ws.Range("B14").PasteSpecial(-4163)
$helper.Clear()
